# Natmi following Dr Hou advice
#
# Re-run of the NATMI ligand-receptor (Tnc -> Ptprb) summary for the
# YoungD0 sample now also includes "ECs" among the sending/target
# clusters (previously only FAPs and sCs were present), per Dr Hou's
# advice. Rewrite sheet1 body (rows 2-10) with the refreshed stats.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
  ,@("ECs", "Tnc", "Ptprb", "ECs", 2, 0.6666666666666666, 1.442371333333333, 4.327114, 0.03522044016446201, 0.03522044016446201, 3, 1, 141.7673063333333, 425.301919, 0.9653621431809045, 0.9653621431809047, 204.4810986590851, 1840.329887931766, 0.03400047960093986, 0.03400047960093987)
  ,@("ECs", "Tnc", "Ptprb", "FAPs", 2, 0.6666666666666666, 1.442371333333333, 4.327114, 0.03522044016446201, 0.03522044016446201, 3, 1, 3.731893333333333, 11.19568, 0.02541226633677049, 0.02541226633677049, 5.382775963057777, 48.44498366752, 0.0008950312059575972, 0.0008950312059575973)
  ,@("ECs", "Tnc", "Ptprb", "sCs", 2, 0.6666666666666666, 1.442371333333333, 4.327114, 0.03522044016446201, 0.03522044016446201, 3, 1, 1.354815, 4.064445, 0.009225590482324892, 0.009225590482324892, 1.95414631797, 17.58731686173, 0.0003249293575645541, 0.0003249293575645541)
  ,@("FAPs", "Tnc", "Ptprb", "ECs", 3, 1, 11.331397, 33.994191, 0.2766948987373093, 0.2766948987373092, 3, 1, 141.7673063333333, 425.301919, 0.9653621431809045, 0.9653621431809047, 1606.421629683614, 14457.79466715253, 0.2671107804522723, 0.2671107804522723)
  ,@("FAPs", "Tnc", "Ptprb", "FAPs", 3, 1, 11.331397, 33.994191, 0.2766948987373093, 0.2766948987373092, 3, 1, 3.731893333333333, 11.19568, 0.02541226633677049, 0.02541226633677049, 42.28756492165333, 380.58808429488, 0.007031444460738244, 0.007031444460738244)
  ,@("FAPs", "Tnc", "Ptprb", "sCs", 3, 1, 11.331397, 33.994191, 0.2766948987373093, 0.2766948987373092, 3, 1, 1.354815, 4.064445, 0.009225590482324892, 0.009225590482324892, 15.351946626555, 138.167519638995, 0.00255267382429877, 0.00255267382429877)
  ,@("sCs", "Tnc", "Ptprb", "ECs", 3, 1, 28.17890933333333, 84.536728, 0.6880846610982287, 0.6880846610982286, 3, 1, 141.7673063333333, 425.301919, 0.9653621431809045, 0.9653621431809047, 3994.848071597892, 35953.63264438103, 0.6642508831276924, 0.6642508831276925)
  ,@("sCs", "Tnc", "Ptprb", "FAPs", 3, 1, 28.17890933333333, 84.536728, 0.6880846610982287, 0.6880846610982286, 3, 1, 3.731893333333333, 11.19568, 0.02541226633677049, 0.02541226633677049, 105.1606838816711, 946.4461549350399, 0.01748579067007465, 0.01748579067007465)
  ,@("sCs", "Tnc", "Ptprb", "sCs", 3, 1, 28.17890933333333, 84.536728, 0.6880846610982287, 0.6880846610982286, 3, 1, 1.354815, 4.064445, 0.009225590482324892, 0.009225590482324892, 38.17720904844001, 343.59488143596, 0.006347987300461568, 0.006347987300461567)
)

for ($i = 0; $i -lt $rows.Length; $i++) {
  $r = $i + 2
  $rowVals = $rows[$i]
  for ($c = 0; $c -lt $rowVals.Length; $c++) {
    $ws.Cells.Item($r, $c + 1).Value = $rowVals[$c]
  }
}
